$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 3 (index 3): reposition / resize the two screenshot pictures
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$s3Pic1 = $s3.Shapes.Item(5)   # "Picture 1"
$s3Pic1.Left = 180.0
$s3Pic1.Width = 222.75009155273438
$s3Pic1.Height = 396.0

$s3Pic2 = $s3.Shapes.Item(6)   # "Picture 2"
$s3Pic2.Left = 504.0
$s3Pic2.Width = 222.75
$s3Pic2.Height = 396.0

# ---------------------------------------------------------------------------
# 2) Slide 4 (index 4): reposition / resize the two screenshot pictures
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$s4Pic1 = $s4.Shapes.Item(4)   # "Picture 1"
$s4Pic1.Left = 180.0
$s4Pic1.Width = 263.2500915527344
$s4Pic1.Height = 468.0

$s4Pic2 = $s4.Shapes.Item(5)   # "Picture 8"
$s4Pic2.Left = 504.0
$s4Pic2.Width = 263.2500915527344
$s4Pic2.Height = 468.0
